# Apply updated numeric data for columns A and B (case 1 dataset refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = -0.1424393201194789
$ws.Cells.Item(1, 2).Value = 0.14230523266662942
$ws.Cells.Item(2, 1).Value = -0.008223192996199558
$ws.Cells.Item(2, 2).Value = 0.008041900043595618
$ws.Cells.Item(3, 1).Value = 0.09269285829147478
$ws.Cells.Item(3, 2).Value = -0.09292453147799407
$ws.Cells.Item(4, 1).Value = -0.19106277818102768
$ws.Cells.Item(4, 2).Value = 0.19017135273812968
$ws.Cells.Item(5, 1).Value = -0.1841713536383054
$ws.Cells.Item(5, 2).Value = 0.18237434989929557
$ws.Cells.Item(6, 1).Value = -0.07951368732317343
$ws.Cells.Item(6, 2).Value = 0.0794352763298738
$ws.Cells.Item(7, 1).Value = -0.05943527741182031
$ws.Cells.Item(7, 2).Value = 0.05927400182362064
$ws.Cells.Item(8, 1).Value = -0.03927400291366112
$ws.Cells.Item(8, 2).Value = 0.03916831614640426
$ws.Cells.Item(9, 1).Value = -0.0331683171005448
$ws.Cells.Item(9, 2).Value = 0.03308869004117554
$ws.Cells.Item(10, 1).Value = -0.02708869100241884
$ws.Cells.Item(10, 2).Value = 0.0270862340331135
$ws.Cells.Item(11, 1).Value = -0.022586234979716124
$ws.Cells.Item(11, 2).Value = 0.022573050991660182
$ws.Cells.Item(12, 1).Value = -0.016573051954973383
$ws.Cells.Item(12, 2).Value = 0.01651803483572456
$ws.Cells.Item(13, 1).Value = -0.01051803580425048
$ws.Cells.Item(13, 2).Value = 0.010499371328323015
$ws.Cells.Item(14, 1).Value = 0.0015006276400830743
$ws.Cells.Item(14, 2).Value = -0.0015317134597587767
$ws.Cells.Item(15, 1).Value = 0.0075317124901408405
$ws.Cells.Item(15, 2).Value = -0.007576942469945003
$ws.Cells.Item(16, 1).Value = -0.015026062731303025
$ws.Cells.Item(16, 2).Value = 0.015003807324497576
$ws.Cells.Item(17, 1).Value = -0.009003808296822235
$ws.Cells.Item(17, 2).Value = 0.00899999899597681
$ws.Cells.Item(18, 1).Value = -0.03610725626583289
$ws.Cells.Item(18, 2).Value = 0.03609577941582387
$ws.Cells.Item(19, 1).Value = -0.02709578033706661
$ws.Cells.Item(19, 2).Value = 0.02701288985222483
$ws.Cells.Item(20, 1).Value = -0.018012890781630375
$ws.Cells.Item(20, 2).Value = 0.01800415533299926
$ws.Cells.Item(21, 1).Value = -0.009004156263469731
$ws.Cells.Item(21, 2).Value = 0.008999999068798559
$ws.Cells.Item(22, 1).Value = -0.09173475925741847
$ws.Cells.Item(22, 2).Value = 0.09147031580245368
$ws.Cells.Item(23, 1).Value = -0.08462223103595345
$ws.Cells.Item(23, 2).Value = 0.0841241262196224
$ws.Cells.Item(24, 1).Value = -0.04212412750568717
$ws.Cells.Item(24, 2).Value = 0.04199999870696214
$ws.Cells.Item(25, 1).Value = -0.096860663510558
$ws.Cells.Item(25, 2).Value = 0.09674339205389515
$ws.Cells.Item(26, 1).Value = -0.09074339299175094
$ws.Cells.Item(26, 2).Value = 0.09059539458135646
$ws.Cells.Item(27, 1).Value = -0.0845953955238139
$ws.Cells.Item(27, 2).Value = 0.08410061264792734
$ws.Cells.Item(28, 1).Value = -0.06608850730254989
$ws.Cells.Item(28, 2).Value = 0.06562659393323145
$ws.Cells.Item(29, 1).Value = -0.05362659496922895
$ws.Cells.Item(29, 2).Value = 0.05341498484868268
$ws.Cells.Item(30, 1).Value = -0.03341498597251524
$ws.Cells.Item(30, 2).Value = 0.0333423710156211
$ws.Cells.Item(31, 1).Value = -0.01834237209455125
$ws.Cells.Item(31, 2).Value = 0.01833578044619877
$ws.Cells.Item(32, 1).Value = -0.006000537061042088
$ws.Cells.Item(32, 2).Value = 0.0059999990130119585

# Narrow column B slightly (closest attainable width to the target 14.7109375 chars)
$ws.Columns.Item(2).ColumnWidth = 13.8
